$wb = $excel.ActiveWorkbook
$wsPrev = $wb.Worksheets.Item("Previously added")
$wsNew = $wb.Worksheets.Item("New")

# ------------------------------------------------------------------
# 1) Archive the 4 listings currently sitting on "New" (rows 2-5) onto
#    the bottom of "Previously added" as rows 335-338, carrying their
#    hyperlinks along.
# ------------------------------------------------------------------

# Clone the formatting of the last existing data row so the appended
# rows look like the rest of the table (style block only - values are
# written separately below).
$wsPrev.Range("A334:F334").Copy()
$wsPrev.Range("A335:F338").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$wsPrev.Cells.Item(335, 1).Value = "https://www.ss.com/msg/lv/real-estate/wood/cesis-and-reg/dzerbenes-pag/ihhmh.html"
$wsPrev.Cells.Item(335, 2).NumberFormat = "@"
$wsPrev.Cells.Item(335, 2).Value = "50 €"
$wsPrev.Cells.Item(335, 3).Value = "Cēsis un raj."
$wsPrev.Cells.Item(335, 4).Value = "10.20 ha."
$wsPrev.Cells.Item(335, 5).NumberFormat = "@"
$wsPrev.Cells.Item(335, 5).Value = "42500040011"
$wsPrev.Cells.Item(335, 6).Value = 46006.506944444445

$wsPrev.Cells.Item(336, 1).Value = "https://www.ss.com/msg/lv/real-estate/wood/cesis-and-reg/vecpiebalgas-pag/alikx.html"
$wsPrev.Cells.Item(336, 2).Value = "21 000 €"
$wsPrev.Cells.Item(336, 3).Value = "Cēsis un raj."
$wsPrev.Cells.Item(336, 4).Value = "3.51 ha."
$wsPrev.Cells.Item(336, 5).NumberFormat = "@"
$wsPrev.Cells.Item(336, 5).Value = "42540020145"
$wsPrev.Cells.Item(336, 6).Value = 46005.57708333334

$wsPrev.Cells.Item(337, 1).Value = "https://www.ss.com/msg/lv/real-estate/wood/gulbene-and-reg/lejasciema-pag/lfkbg.html"
$wsPrev.Cells.Item(337, 2).Value = "26 400 €"
$wsPrev.Cells.Item(337, 3).Value = "Gulbene un raj."
$wsPrev.Cells.Item(337, 4).Value = "7 ha."
$wsPrev.Cells.Item(337, 5).NumberFormat = "@"
$wsPrev.Cells.Item(337, 5).Value = "50640020038"
$wsPrev.Cells.Item(337, 6).Value = 46003.92222222222

$wsPrev.Cells.Item(338, 1).Value = "https://www.ss.com/msg/lv/real-estate/wood/ludza-and-reg/ciblas-pag/ggibg.html"
$wsPrev.Cells.Item(338, 2).Value = "85 000 €"
$wsPrev.Cells.Item(338, 3).Value = "Ludza un raj."
$wsPrev.Cells.Item(338, 4).Value = "15.70 ha."
$wsPrev.Cells.Item(338, 5).NumberFormat = "@"
$wsPrev.Cells.Item(338, 5).Value = "68480030023"
$wsPrev.Cells.Item(338, 6).Value = 46003.8625

$wsPrev.Hyperlinks.Add($wsPrev.Cells.Item(335, 1), "https://www.ss.com/msg/lv/real-estate/wood/cesis-and-reg/dzerbenes-pag/ihhmh.html")
$wsPrev.Hyperlinks.Add($wsPrev.Cells.Item(336, 1), "https://www.ss.com/msg/lv/real-estate/wood/cesis-and-reg/vecpiebalgas-pag/alikx.html")
$wsPrev.Hyperlinks.Add($wsPrev.Cells.Item(337, 1), "https://www.ss.com/msg/lv/real-estate/wood/gulbene-and-reg/lejasciema-pag/lfkbg.html")
$wsPrev.Hyperlinks.Add($wsPrev.Cells.Item(338, 1), "https://www.ss.com/msg/lv/real-estate/wood/ludza-and-reg/ciblas-pag/ggibg.html")

# Adding the hyperlinks re-styles column A with the generic "Hyperlink"
# theme style, so restamp the whole block's formatting once more,
# afterwards, to get back the workbook's own custom look.
$wsPrev.Range("A334:F334").Copy()
$wsPrev.Range("A335:F338").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ------------------------------------------------------------------
# 2) Rebuild "New": drop the old hyperlinks/rows, keep only the single
#    freshly scraped listing on row 2.
# ------------------------------------------------------------------
$wsNew.Hyperlinks.Delete()
$wsNew.Rows("3:5").Delete()

$wsNew.Cells.Item(2, 1).Value = "https://www.ss.com/msg/lv/real-estate/wood/rezekne-and-reg/kaunatas-pag/nnplp.html"
$wsNew.Cells.Item(2, 2).Value = "25 000 €"
$wsNew.Cells.Item(2, 3).Value = "Rēzekne un raj."
$wsNew.Cells.Item(2, 4).Value = "4 ha."
$wsNew.Cells.Item(2, 5).NumberFormat = "@"
$wsNew.Cells.Item(2, 5).Value = "78620020231"
$wsNew.Cells.Item(2, 6).Value = 46006.62986111111

$wsNew.Hyperlinks.Add($wsNew.Cells.Item(2, 1), "https://www.ss.com/msg/lv/real-estate/wood/rezekne-and-reg/kaunatas-pag/nnplp.html")

# Restore the custom per-column formatting on row 2 (same fixup as above).
$wsPrev.Range("A334:F334").Copy()
$wsNew.Range("A2:F2").PasteSpecial(-4122)
$excel.CutCopyMode = $false
